# The deck's single live theme (backing ppt/theme/theme2.xml, the one actually
# wired to the slide master / presentation) is the "Integral" green palette.
# A second, orphaned theme part (ppt/theme/theme1.xml, only referenced by the
# notes master) holds the stock "Office Theme" blue/orange palette. The
# authored change swaps the two parts' contents wholesale (file names stay
# put; the palettes trade places) - font scheme and format scheme are
# byte-identical between the two parts already, so the only observable effect
# is the 12-slot color scheme driving the live theme flipping from the
# Integral colors to the stock Office Theme colors.
#
# Reproduce that through the object model by rewriting the live theme's
# ThemeColorScheme in place with the stock Office Theme RGB values (colors
# are COM RGB() integers, i.e. 0x00BBGGRR, decimal-encoded).

$p = $ppt.ActivePresentation
$design = $p.Designs.Item(1)
$master = $design.SlideMaster
$theme = $master.Theme
$colors = $theme.ThemeColorScheme

$colors.Item(1).RGB = 0           # dk1      000000
$colors.Item(2).RGB = 16777215    # lt1      FFFFFF
$colors.Item(3).RGB = 6968388     # dk2      44546A
$colors.Item(4).RGB = 15132391    # lt2      E7E6E6
$colors.Item(5).RGB = 13998939    # accent1  5B9BD5
$colors.Item(6).RGB = 3243501     # accent2  ED7D31
$colors.Item(7).RGB = 10855845    # accent3  A5A5A5
$colors.Item(8).RGB = 49407       # accent4  FFC000
$colors.Item(9).RGB = 12874308    # accent5  4472C4
$colors.Item(10).RGB = 4697456    # accent6  70AD47
$colors.Item(11).RGB = 12673797   # hlink    0563C1
$colors.Item(12).RGB = 7491477    # folHlink 954F72
